$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.080.73'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.98%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.832.14'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.38%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.38'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -4.51%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4594'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -4.78%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2674'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -6.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06172'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -5.80%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.833.25'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -4.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07337'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.85%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -4.96%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -4.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '82.23'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -6.66%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6151'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -7.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.039.04'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.001'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '224.66'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.93%  '
$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'BinanceUSD'
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.002'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.081.43'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -6.71%  '
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000007175'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -5.50%  '
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.25'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -8.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.809'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -9.04%  '
$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.812'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -6.68%  '
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '164.89'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.05%  '
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.069'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.00%  '
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.52'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -6.71%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.832'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -6.86%  '
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'Stellar'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1008'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.71%  '
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.371'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.17%  '
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.047'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -7.07%  '
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.741'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -7.27%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04777'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -6.50%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.128'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -7.07%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6896'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -9.13%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.001'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.19%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.696'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.62%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01796'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -4.39%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.608'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.62%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8814'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -4.29%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.900'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -8.57%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9960'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.85%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '102.71'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -4.34%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.410'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -5.16%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3965'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -7.87%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.830'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -8.11%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1174'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -7.91%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '58.79'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -9.06%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.425'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -6.25%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05531'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.55%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '32.34'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -4.73%  '
